$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string need NumberFormat forced
# to Text first, otherwise Excel auto-converts them to a numeric value.
$textCells = @('D5', 'D6', 'D8', 'D10', 'D12', 'D14', 'D17', 'D18', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D31', 'D32', 'D34', 'D35', 'D36', 'D37', 'D42', 'D43', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '72.136.69'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '4.039.17'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '539.48'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '152.02'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('D7').Value = '4.033.20'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.699'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = '0.754'
$ws.Range('E10').Value = '  -1.39%  '
$ws.Range('E11').Value = '  -1.34%  '
$ws.Range('D12').Value = '54.13'
$ws.Range('E12').Value = '  +11.40%  '
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = '10.91'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').Value = '4.685.42'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '4.044.43'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '14.35'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '20.59'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D21').Value = '72.133.03'
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('D22').Value = '444.89'
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('D23').Value = '97.27'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -2.29%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '4.25'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '14.68'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').Value = '4.30'
$ws.Range('E27').Value = '  +16.28%  '
$ws.Range('D28').Value = '11.28'
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '10.78'
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('E30').Value = '  +1.68%  '
$ws.Range('D31').Value = '37.15'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '8.19'
$ws.Range('E32').Value = '  +18.12%  '
$ws.Range('E33').Value = '  +1.19%  '
$ws.Range('D34').Value = '13.56'
$ws.Range('E34').Value = '  -1.34%  '
$ws.Range('D35').Value = '49.28'
$ws.Range('E35').Value = '  +14.24%  '
$ws.Range('D36').Value = '681.21'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '66.88'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('E38').Value = '  +4.44%  '
$ws.Range('D39').Value = '0.0₃0874'
$ws.Range('E39').Value = '  +2.72%  '
$ws.Range('E40').Value = '  -5.69%  '
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '3.38'
$ws.Range('E42').Value = '  -3.98%  '
$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D43').Value = '11.23'
$ws.Range('E43').Value = '  +17.49%  '
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('E48').Value = '  -3.77%  '
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('D50').Value = '3.31'
$ws.Range('E50').Value = '  -3.20%  '
$ws.Range('E51').Value = '  +1.03%  '
